$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# GDP per Capita values for years 1950-2016 (row 2 = 1950 ... row 68 = 2016).
# The source data changed for most years and six new years (2011-2016) were
# appended to the bottom of the table.
$gdpValues = @("676","743","740","697","735","733","724","720","708","724","732","703","746","770","787","792","867","875","894","883","905","915","945","945","937","956","996","993","993","977","958","920","905","875","870","843","845","859","867","870","875","891.894264796583","895.480598700604","900.413342270506","904.757464651378","931.002765921089","979.730435078623","1017.16811846792","1052.07501679565","1100.35320685852","1148.68017009953","1213.62687531406","1293.20506586684","1369.93908161501","1464.38536693869","1554.15231589586","1621.99654449678","1754.462648525","1848.17452737512","1946.84813162394","2068.86366158353","2227","2275","2372","2467","2566","2660")

$firstRow = 2
$firstYear = 1950
$lastRow = $firstRow + $gdpValues.Length - 1

# The "Data" column holds numeric-looking text (shared strings in the original
# file), so mark the range as text before assigning values to avoid Excel
# auto-converting them to numbers, then drop the temporary formatting again.
$dataRange = $ws.Range("E" + $firstRow + ":E" + $lastRow)
$dataRange.NumberFormat = "@"

for ($i = 0; $i -lt $gdpValues.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("A$row").Value = 834
    $ws.Range("B$row").Value = "Tanzania"
    $ws.Range("C$row").Value = "GDP per Capita"
    $ws.Range("D$row").Value = $firstYear + $i
    $ws.Range("E$row").Value = $gdpValues[$i]
}

$dataRange.ClearFormats()
